# Natmi following Dr Hou advice
#
# The LR-pair result table (Sending cluster / Ligand / Receptor / Target
# cluster + the NATMI specificity metrics) is recomputed for a 3-cluster
# model (ECs, FAPs, sCs) instead of the previous 2-cluster one, which
# turns the old 2-row (A1:T3) edge list into a full 3x3 = 9-row
# (A1:T10) edge list for the Fgf18 -> Fgfr4 pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the NATMI numeric metrics (unchanged headers
# in row 1).
$data = @(
    @("ECs","Fgf18","Fgfr4","ECs",  1,0.3333333333333333,0.4765323333333333,1.429597,          0.04733670697480491,0.04733670697480491,1,0.3333333333333333,0.132884,          0.398652,          0.01195569974366677,0.01195569974366677,0.06332352258266667,0.569911703244,    0.0005659434554447038,0.0005659434554447038),
    @("ECs","Fgf18","Fgfr4","FAPs", 1,0.3333333333333333,0.4765323333333333,1.429597,          0.04733670697480491,0.04733670697480491,1,0.3333333333333333,0.07352966666666667,0.220589,          0.006615533976389704,0.006615533976389703,0.03503926362588889,0.315353372633,    0.0003131575933222254,0.0003131575933222253),
    @("ECs","Fgf18","Fgfr4","sCs",  1,0.3333333333333333,0.4765323333333333,1.429597,          0.04733670697480491,0.04733670697480491,3,1,                  10.908285,         32.724855,         0.9814287662799436, 0.9814287662799435, 5.198150503715,    46.78335453343499,  0.04645760592603798,  0.04645760592603797),
    @("FAPs","Fgf18","Fgfr4","ECs", 3,1,                  8.770835333333332,  26.312506,         0.8712576945074703, 0.8712576945074703, 1,0.3333333333333333,0.132884,          0.398652,          0.01195569974366677,0.01195569974366677,1.165503682434667, 10.489533141912,    0.01041649539489066,  0.01041649539489066),
    @("FAPs","Fgf18","Fgfr4","FAPs",3,1,                  8.770835333333332,  26.312506,         0.8712576945074703, 0.8712576945074703, 1,0.3333333333333333,0.07352966666666667,0.220589,          0.006615533976389704,0.006615533976389703,0.6449165984482222,5.804249386034,    0.005763834880205131, 0.00576383488020513),
    @("FAPs","Fgf18","Fgfr4","sCs", 3,1,                  8.770835333333332,  26.312506,         0.8712576945074703, 0.8712576945074703, 3,1,                  10.908285,         32.724855,         0.9814287662799436, 0.9814287662799435, 95.67477150406998, 861.07294353663,    0.8550773642323745,   0.8550773642323745),
    @("sCs","Fgf18","Fgfr4","ECs",  3,1,                  0.8194993333333332, 2.458498,          0.08140559851772487,0.08140559851772486,1,0.3333333333333333,0.132884,          0.398652,          0.01195569974366677,0.01195569974366677,0.1088983494106667,0.9800851446959998, 0.0009732608933314029,0.0009732608933314028),
    @("sCs","Fgf18","Fgfr4","FAPs", 3,1,                  0.8194993333333332, 2.458498,          0.08140559851772487,0.08140559851772486,1,0.3333333333333333,0.07352966666666667,0.220589,          0.006615533976389704,0.006615533976389703,0.06025751281355556,0.542317615322,    0.0005385415028623482,0.000538541502862348),
    @("sCs","Fgf18","Fgfr4","sCs",  3,1,                  0.8194993333333332, 2.458498,          0.08140559851772487,0.08140559851772486,3,1,                  10.908285,         32.724855,         0.9814287662799436, 0.9814287662799435, 8.939332285309998, 80.45399056778999,  0.07989379612153112,  0.07989379612153111)
)

# Rows 2 and 3 already existed (ECs-vs-ECs, ECs-vs-FAPs under their old
# FAPs/sCs cluster names); overwrite them in place and append the new
# rows 4-10 so the sheet grows from A1:T3 to A1:T10.
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
